$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B30 was stored as a text string "5"; convert it to a true numeric value 5.
$ws.Range("B30").Value = 5

# Add a new row 31 with the new annotation data.
$ws.Range("A31").Value = "Sunsi Wu"

# B31 keeps "3" as text (matches source data format), so force a text
# number format before assignment to avoid Excel's automatic numeric
# coercion, then restore the default style so no stray formatting sticks.
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "3"
$ws.Range("B31").Style = "Normal"

$ws.Range("C31").Value = "are not properly"
$ws.Range("D31").Value = "DFT"
$ws.Range("E31").Value = "WRI"
$ws.Range("F31").Value = "74483628-1e12-4bb7-acfc-2ccaf38e6d81"
$ws.Range("G31").Value = "HyIFzx-0b_annotated.xlsx"
$ws.Range("H31").Value = "Acronyms are not properly defined."
